# Add newly uploaded product rows to the "SIMAR" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SIMAR")
$ws.Activate()

# New product data: Name (A), Category (B), Price (C), Description (D), Image URL (E)
$newRows = @(
    @{ Name = "J SONS/A-2528";   Category = "TEXTILE SUIT"; Price = 710; Description = "J SONS/A-2528 - 16/20 - 3PCS BOX - IMPORTED - 710";   Image = "https://i.postimg.cc/zX92vhR0/Whats-App-Image-2025-05-30-at-11-04-56.jpg" },
    @{ Name = "CRISPY/A-30859";  Category = "TEXTILE SUIT"; Price = 785; Description = "CRISPY/A-30859 - 22/32 - 1PCS BOX - 785";              Image = "https://i.postimg.cc/mknXQzpZ/Whats-App-Image-2025-05-30-at-11-04-56-1.jpg" },
    @{ Name = "J SONS/A-2530";   Category = "TEXTILE SUIT"; Price = 755; Description = "J SONS/A-2530 - 16/20 - 1PCS BOX - IMPORTED - 755";   Image = "https://i.postimg.cc/pVc6N4VN/Whats-App-Image-2025-05-30-at-11-04-57.jpg" },
    @{ Name = "J SONS/A-3438";   Category = "TEXTILE SUIT"; Price = 695; Description = "J SONS/A-3438 - 22/32 - 3PCS BOX - IMPORTED - 695";   Image = "https://i.postimg.cc/kXVh7sD2/Whats-App-Image-2025-05-30-at-11-04-58.jpg" }
)

$startRow = 126
$r = $startRow
foreach ($row in $newRows) {
    $ws.Range("A$r").Value = $row.Name
    $ws.Range("B$r").Value = $row.Category
    $ws.Range("C$r").Value = $row.Price
    $ws.Range("D$r").Value = $row.Description
    $ws.Range("E$r").Value = $row.Image
    $r = $r + 1
}

$lastRow = $r - 1

# Reflect the sheet view state from the saved workbook (best effort).
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 112
$win.ScrollColumn = 1
$ws.Range("D120").Select()
